$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64, shifting existing rows 64..145 down to 65..146.
$ws.Rows("64:64").Insert()

# Populate the newly inserted row 64 with the new weekly price observation.
$ws.Cells.Item(64, 1).Value = 8
$ws.Cells.Item(64, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44915
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 1200
$ws.Cells.Item(64, 11).Value = 2500
$ws.Cells.Item(64, 12).Value = 3000
$ws.Cells.Item(64, 13).Value = 2750
$ws.Cells.Item(64, 14).Value = "`$/paquete"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 2750
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
